$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New public-exposure-site rows replacing the old Mount Waverley entries.
$data = @(
    @("Black Rock",      "Woolworths Metro  40 Bluff Road, Black Rock VIC 3193",       "30/12/20 5:30pm-5:55pm",     "Case shopped",          "new"),
    @("Box Hill South",  "Bunnings  259 Middleborough Road, Box Hill South VIC 3128",  "30/12/20 12:00pm-12:40pm",   "Case shopped",          "new"),
    @("Cheltenham",      "Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192",     "29/12/20 01:30pm-01:45pm",   "Case shopped in store", "new"),
    @("Cheltenham",      "Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192",     "29/12/2020 01:30pm-01:45pm", "Case shopped in store", "old")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

# Match the column widths that Excel's "best fit" recalculated for the new,
# longer location / site / exposure-period / notes text.
$ws.Columns.Item(1).ColumnWidth = 11.73046875
$ws.Columns.Item(2).ColumnWidth = 48.9296875
$ws.Columns.Item(3).ColumnWidth = 25.6640625
$ws.Columns.Item(4).ColumnWidth = 18

$ws.Range("B5").Select()
